$wb = $excel.ActiveWorkbook

# --- LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(39, 1).Value = "06:21:22"
$ws.Cells.Item(39, 2).Value = "06:29"
$ws.Cells.Item(39, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(39, 4).Value = 8
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(40, 1).Value = "06:21:22"
$ws.Cells.Item(40, 2).Value = "06:29"
$ws.Cells.Item(40, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(40, 4).Value = 8
$ws.Cells.Item(40, 5).Value = "LP1912"
$ws.Cells.Item(47, 1).Value = "06:59:37"
$ws.Cells.Item(47, 2).Value = "06:59"
$ws.Cells.Item(47, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = "LP1912"
$ws.Cells.Item(48, 1).Value = "06:59:37"
$ws.Cells.Item(48, 2).Value = "06:59"
$ws.Cells.Item(48, 3).Value = "14_ABASTO"
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = "LP1912"
$ws.Cells.Item(63, 1).Value = "06:59:37"
$ws.Cells.Item(63, 2).Value = "07:31"
$ws.Cells.Item(63, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(63, 4).Value = 32
$ws.Cells.Item(63, 5).Value = "LP1912"
$ws.Cells.Item(64, 1).Value = "06:59:37"
$ws.Cells.Item(64, 2).Value = "07:31"
$ws.Cells.Item(64, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(64, 4).Value = 32
$ws.Cells.Item(64, 5).Value = "LP1912"
$ws.Cells.Item(65, 1).Value = "05:52:07"
$ws.Cells.Item(65, 2).Value = "07:32"
$ws.Cells.Item(65, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(65, 4).Value = 100
$ws.Cells.Item(65, 5).Value = "LP1912"
$ws.Cells.Item(67, 1).Value = "05:52:07"
$ws.Cells.Item(67, 2).Value = "07:32"
$ws.Cells.Item(67, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(67, 4).Value = 100
$ws.Cells.Item(67, 5).Value = "LP1912"
$ws.Cells.Item(75, 1).Value = "07:51:34"
$ws.Cells.Item(75, 2).Value = "07:51"
$ws.Cells.Item(75, 3).Value = "10_OLMOS"
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = "LP1912"
$ws.Cells.Item(76, 1).Value = "07:51:34"
$ws.Cells.Item(76, 2).Value = "07:51"
$ws.Cells.Item(76, 3).Value = "215D_EL PATO"
$ws.Cells.Item(76, 4).Value = 0
$ws.Cells.Item(76, 5).Value = "LP1912"
$ws.Cells.Item(86, 1).Value = "08:13:38"
$ws.Cells.Item(86, 2).Value = "08:21"
$ws.Cells.Item(86, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(86, 4).Value = 8
$ws.Cells.Item(86, 5).Value = "LP1912"
$ws.Cells.Item(87, 1).Value = "08:13:38"
$ws.Cells.Item(87, 2).Value = "08:22"
$ws.Cells.Item(87, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(87, 4).Value = 9
$ws.Cells.Item(87, 5).Value = "LP1912"
$ws.Cells.Item(88, 1).Value = "08:13:38"
$ws.Cells.Item(88, 2).Value = "08:23"
$ws.Cells.Item(88, 3).Value = "215B_EL PATO"
$ws.Cells.Item(88, 4).Value = 10
$ws.Cells.Item(88, 5).Value = "LP1912"
$ws.Cells.Item(90, 1).Value = "08:13:38"
$ws.Cells.Item(90, 2).Value = "08:27"
$ws.Cells.Item(90, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(90, 4).Value = 14
$ws.Cells.Item(90, 5).Value = "LP1912"
$ws.Cells.Item(92, 1).Value = "08:13:38"
$ws.Cells.Item(92, 2).Value = "08:33"
$ws.Cells.Item(92, 3).Value = "10_OLMOS"
$ws.Cells.Item(92, 4).Value = 20
$ws.Cells.Item(92, 5).Value = "LP1912"
$ws.Cells.Item(93, 1).Value = "08:13:38"
$ws.Cells.Item(93, 2).Value = "08:36"
$ws.Cells.Item(93, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(93, 4).Value = 23
$ws.Cells.Item(93, 5).Value = "LP1912"
$ws.Cells.Item(94, 1).Value = "08:13:38"
$ws.Cells.Item(94, 2).Value = "08:42"
$ws.Cells.Item(94, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(94, 4).Value = 29
$ws.Cells.Item(94, 5).Value = "LP1912"
$ws.Cells.Item(95, 1).Value = "08:13:38"
$ws.Cells.Item(95, 2).Value = "08:43"
$ws.Cells.Item(95, 3).Value = "14_ABASTO"
$ws.Cells.Item(95, 4).Value = 30
$ws.Cells.Item(95, 5).Value = "LP1912"
$ws.Cells.Item(96, 1).Value = "07:28:14"
$ws.Cells.Item(96, 2).Value = "08:44"
$ws.Cells.Item(96, 3).Value = "14_ABASTO"
$ws.Cells.Item(96, 4).Value = 76
$ws.Cells.Item(96, 5).Value = "LP1912"
$ws.Cells.Item(97, 1).Value = "08:13:38"
$ws.Cells.Item(97, 2).Value = "08:53"
$ws.Cells.Item(97, 3).Value = "10_OLMOS"
$ws.Cells.Item(97, 4).Value = 40
$ws.Cells.Item(97, 5).Value = "LP1912"
$ws.Cells.Item(98, 1).Value = "08:13:38"
$ws.Cells.Item(98, 2).Value = "08:54"
$ws.Cells.Item(98, 3).Value = "17_ROMERO"
$ws.Cells.Item(98, 4).Value = 41
$ws.Cells.Item(98, 5).Value = "LP1912"
$ws.Cells.Item(99, 1).Value = "08:13:38"
$ws.Cells.Item(99, 2).Value = "09:01"
$ws.Cells.Item(99, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99, 4).Value = 48
$ws.Cells.Item(99, 5).Value = "LP1912"
$ws.Cells.Item(100, 1).Value = "08:13:38"
$ws.Cells.Item(100, 2).Value = "09:01"
$ws.Cells.Item(100, 3).Value = "215A_EL PATO"
$ws.Cells.Item(100, 4).Value = 48
$ws.Cells.Item(100, 5).Value = "LP1912"
$ws.Cells.Item(101, 1).Value = "07:28:14"
$ws.Cells.Item(101, 2).Value = "09:02"
$ws.Cells.Item(101, 3).Value = "215A_EL PATO"
$ws.Cells.Item(101, 4).Value = 94
$ws.Cells.Item(101, 5).Value = "LP1912"
$ws.Cells.Item(102, 1).Value = "08:13:38"
$ws.Cells.Item(102, 2).Value = "09:03"
$ws.Cells.Item(102, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(102, 4).Value = 50
$ws.Cells.Item(102, 5).Value = "LP1912"
$ws.Cells.Item(103, 1).Value = "08:13:38"
$ws.Cells.Item(103, 2).Value = "09:10"
$ws.Cells.Item(103, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(103, 4).Value = 57
$ws.Cells.Item(103, 5).Value = "LP1912"
$ws.Cells.Item(104, 1).Value = "07:28:14"
$ws.Cells.Item(104, 2).Value = "09:11"
$ws.Cells.Item(104, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(104, 4).Value = 103
$ws.Cells.Item(104, 5).Value = "LP1912"
$ws.Cells.Item(105, 1).Value = "08:13:38"
$ws.Cells.Item(105, 2).Value = "09:16"
$ws.Cells.Item(105, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(105, 4).Value = 63
$ws.Cells.Item(105, 5).Value = "LP1912"
$ws.Cells.Item(106, 1).Value = "07:28:14"
$ws.Cells.Item(106, 2).Value = "09:17"
$ws.Cells.Item(106, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(106, 4).Value = 109
$ws.Cells.Item(106, 5).Value = "LP1912"
$ws.Cells.Item(107, 1).Value = "08:13:38"
$ws.Cells.Item(107, 2).Value = "09:21"
$ws.Cells.Item(107, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(107, 4).Value = 68
$ws.Cells.Item(107, 5).Value = "LP1912"
$ws.Cells.Item(108, 1).Value = "08:13:38"
$ws.Cells.Item(108, 2).Value = "09:22"
$ws.Cells.Item(108, 3).Value = "17_ROMERO"
$ws.Cells.Item(108, 4).Value = 69
$ws.Cells.Item(108, 5).Value = "LP1912"
$ws.Cells.Item(109, 1).Value = "07:28:14"
$ws.Cells.Item(109, 2).Value = "09:23"
$ws.Cells.Item(109, 3).Value = "17_ROMERO"
$ws.Cells.Item(109, 4).Value = 115
$ws.Cells.Item(109, 5).Value = "LP1912"
$ws.Cells.Item(110, 1).Value = "08:13:38"
$ws.Cells.Item(110, 2).Value = "09:23"
$ws.Cells.Item(110, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(110, 4).Value = 70
$ws.Cells.Item(110, 5).Value = "LP1912"
$ws.Cells.Item(111, 1).Value = "07:28:14"
$ws.Cells.Item(111, 2).Value = "09:24"
$ws.Cells.Item(111, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(111, 4).Value = 116
$ws.Cells.Item(111, 5).Value = "LP1912"
$ws.Cells.Item(112, 1).Value = "08:13:38"
$ws.Cells.Item(112, 2).Value = "09:32"
$ws.Cells.Item(112, 3).Value = "15_ABASTO"
$ws.Cells.Item(112, 4).Value = 79
$ws.Cells.Item(112, 5).Value = "LP1912"
$ws.Cells.Item(113, 1).Value = "08:13:38"
$ws.Cells.Item(113, 2).Value = "09:33"
$ws.Cells.Item(113, 3).Value = "10_OLMOS"
$ws.Cells.Item(113, 4).Value = 80
$ws.Cells.Item(113, 5).Value = "LP1912"
$ws.Cells.Item(114, 1).Value = "08:13:38"
$ws.Cells.Item(114, 2).Value = "09:34"
$ws.Cells.Item(114, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(114, 4).Value = 81
$ws.Cells.Item(114, 5).Value = "LP1912"
$ws.Cells.Item(115, 1).Value = "08:13:38"
$ws.Cells.Item(115, 2).Value = "09:41"
$ws.Cells.Item(115, 3).Value = "215C_EL PATO"
$ws.Cells.Item(115, 4).Value = 88
$ws.Cells.Item(115, 5).Value = "LP1912"
$ws.Cells.Item(116, 1).Value = "07:51:34"
$ws.Cells.Item(116, 2).Value = "09:42"
$ws.Cells.Item(116, 3).Value = "215C_EL PATO"
$ws.Cells.Item(116, 4).Value = 111
$ws.Cells.Item(116, 5).Value = "LP1912"
$ws.Cells.Item(117, 1).Value = "08:13:38"
$ws.Cells.Item(117, 2).Value = "09:43"
$ws.Cells.Item(117, 3).Value = "14_ABASTO"
$ws.Cells.Item(117, 4).Value = 90
$ws.Cells.Item(117, 5).Value = "LP1912"
$ws.Cells.Item(118, 1).Value = "08:13:38"
$ws.Cells.Item(118, 2).Value = "09:58"
$ws.Cells.Item(118, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(118, 4).Value = 105
$ws.Cells.Item(118, 5).Value = "LP1912"
$ws.Cells.Item(119, 1).Value = "08:13:38"
$ws.Cells.Item(119, 2).Value = "10:12"
$ws.Cells.Item(119, 3).Value = "15_ABASTO"
$ws.Cells.Item(119, 4).Value = 119
$ws.Cells.Item(119, 5).Value = "LP1912"

# --- LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(21, 1).Value = "08:13:38"
$ws.Cells.Item(21, 2).Value = "08:23"
$ws.Cells.Item(21, 3).Value = "215B_EL PATO"
$ws.Cells.Item(21, 4).Value = 10
$ws.Cells.Item(21, 5).Value = "LP1912"
$ws.Cells.Item(22, 1).Value = "08:13:38"
$ws.Cells.Item(22, 2).Value = "09:01"
$ws.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws.Cells.Item(22, 4).Value = 48
$ws.Cells.Item(22, 5).Value = "LP1912"
$ws.Cells.Item(24, 1).Value = "08:13:38"
$ws.Cells.Item(24, 2).Value = "09:41"
$ws.Cells.Item(24, 3).Value = "215C_EL PATO"
$ws.Cells.Item(24, 4).Value = 88
$ws.Cells.Item(24, 5).Value = "LP1912"
$ws.Cells.Item(25, 1).Value = "07:51:34"
$ws.Cells.Item(25, 2).Value = "09:42"
$ws.Cells.Item(25, 3).Value = "215C_EL PATO"
$ws.Cells.Item(25, 4).Value = 111
$ws.Cells.Item(25, 5).Value = "LP1912"

# --- 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(21, 1).Value = "08:13:38"
$ws.Cells.Item(21, 2).Value = "08:16"
$ws.Cells.Item(21, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(21, 4).Value = 3
$ws.Cells.Item(21, 5).Value = "L6203"
$ws.Cells.Item(22, 1).Value = "06:49:33"
$ws.Cells.Item(22, 2).Value = "08:33"
$ws.Cells.Item(22, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(22, 4).Value = 104
$ws.Cells.Item(22, 5).Value = "L6173"
$ws.Cells.Item(23, 1).Value = "07:28:14"
$ws.Cells.Item(23, 2).Value = "08:38"
$ws.Cells.Item(23, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(23, 4).Value = 70
$ws.Cells.Item(23, 5).Value = "L6173"
$ws.Cells.Item(24, 1).Value = "07:51:34"
$ws.Cells.Item(24, 2).Value = "08:40"
$ws.Cells.Item(24, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(24, 4).Value = 49
$ws.Cells.Item(24, 5).Value = "L6173"
$ws.Cells.Item(25, 1).Value = "08:13:38"
$ws.Cells.Item(25, 2).Value = "08:45"
$ws.Cells.Item(25, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(25, 4).Value = 32
$ws.Cells.Item(25, 5).Value = "L6173"
$ws.Cells.Item(26, 1).Value = "08:13:38"
$ws.Cells.Item(26, 2).Value = "09:08"
$ws.Cells.Item(26, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(26, 4).Value = 55
$ws.Cells.Item(26, 5).Value = "L6203"
$ws.Cells.Item(27, 1).Value = "07:51:34"
$ws.Cells.Item(27, 2).Value = "09:09"
$ws.Cells.Item(27, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(27, 4).Value = 78
$ws.Cells.Item(27, 5).Value = "L6203"
$ws.Cells.Item(28, 1).Value = "08:13:38"
$ws.Cells.Item(28, 2).Value = "10:03"
$ws.Cells.Item(28, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(28, 4).Value = 110
$ws.Cells.Item(28, 5).Value = "L6173"

# --- header rows (Última actualización / Total filas) ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:13:38"
$ws.Cells.Item(3, 1).Value = "Total filas: 114"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:13:38"
$ws.Cells.Item(3, 1).Value = "Total filas: 20"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:13:38"
$ws.Cells.Item(3, 1).Value = "Total filas: 23"
